$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.169.85"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "3.854.99"
$ws.Range("E3").Value = "  +1.10%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "697.76"
$ws.Range("E5").Value = "  +2.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.28"
$ws.Range("E6").Value = "  +0.99%  "

$ws.Range("D7").Value = "3.853.47"
$ws.Range("E7").Value = "  +1.10%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.14"
$ws.Range("E11").Value = "  -1.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("E13").Value = "  +4.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.38"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").Value = "4.506.20"
$ws.Range("E15").Value = "  +1.13%  "

$ws.Range("D16").Value = "3.852.30"
$ws.Range("E16").Value = "  +1.05%  "

$ws.Range("D17").Value = "71.292.65"
$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.13"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "491.83"
$ws.Range("E22").Value = "  +2.92%  "

$ws.Range("E23").Value = "  +1.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.16"
$ws.Range("E24").Value = "  +2.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000146"
$ws.Range("E25").Value = "  +1.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.30"
$ws.Range("E26").Value = "  +0.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.68"
$ws.Range("E27").Value = "  +3.34%  "

$ws.Range("E28").Value = "  +1.52%  "

$ws.Range("D29").Value = "4.012.52"
$ws.Range("E29").Value = "  +1.21%  "

$ws.Range("E30").Value = "  +9.51%  "

$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.65"
$ws.Range("E32").Value = "  +3.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.28"
$ws.Range("E33").Value = "  -0.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.68"
$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.178"
$ws.Range("E35").Value = "  -0.76%  "

$ws.Range("E36").Value = "  +1.58%  "

$ws.Range("D37").Value = "3.807.70"
$ws.Range("E37").Value = "  +1.00%  "

$ws.Range("E38").Value = "  +0.76%  "

$ws.Range("E39").Value = "  +1.70%  "

$ws.Range("E40").Value = "  +11.48%  "

$ws.Range("E41").Value = "  +1.78%  "

$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("E43").Value = "  +6.88%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.57"
$ws.Range("E46").Value = "  +2.25%  "

$ws.Range("E47").Value = "  +3.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.68"
$ws.Range("E48").Value = "  +1.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.33"
$ws.Range("E49").Value = "  -3.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "419.97"
$ws.Range("E50").Value = "  +5.41%  "

$ws.Range("E51").Value = "  +1.06%  "
